# The document is a title-page layout, so each footer/header pair has a
# "first page" part (footer1.xml / header1.xml) distinct from the
# "default" part used on the rest of the pages (footer2.xml / header2.xml).
# The commit renames the InlineShape's display Name on three pictures:
#   footer "first"    (footer1.xml): Pearson logo  image1.png -> image2.png
#   footer "default"  (footer2.xml): Pearson logo  image1.png -> image2.png
#   header "first"    (header1.xml): BTec logo     image2.jpg -> image1.jpg
#
# Renaming an InlineShape that lives in a footer/header story only commits
# reliably once the shape's Range has been Select()-ed and the rename is
# then applied through $word.Selection.InlineShapes - setting .Name
# directly on the shape object fetched straight from the
# Header/Footer.Range.InlineShapes collection does not persist for these
# footer/header stories, so route every rename through Selection.

function Rename-InlineLogo($range, $newName) {
    $shape = $range.InlineShapes.Item(1)
    $shape.Range.Select()
    $word.Selection.InlineShapes.Item(1).Name = $newName
}

$d = $word.ActiveDocument
$section = $d.Sections(1)

# wdHeaderFooterPrimary = 1, wdHeaderFooterFirstPage = 2
Rename-InlineLogo $section.Footers.Item(2).Range "image2.png"
Rename-InlineLogo $section.Footers.Item(1).Range "image2.png"
Rename-InlineLogo $section.Headers.Item(2).Range "image1.jpg"
